$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 13).Value = 1.11  # M2: 1.1 -> 1.11
$ws.Cells.Item(2, 14).Value = 6.5  # N2: 7 -> 6.5
$ws.Cells.Item(4, 7).Value = 4.2  # G4: 4.1 -> 4.2
$ws.Cells.Item(4, 9).Value = 2  # I4: 2.05 -> 2
$ws.Cells.Item(4, 32).Value = 81  # AF4: 67 -> 81
$ws.Cells.Item(4, 35).Value = 8  # AI4: 8.5 -> 8
$ws.Cells.Item(4, 37).Value = 17  # AK4: 19 -> 17
$ws.Cells.Item(4, 44).Value = 151  # AR4: 126 -> 151
$ws.Cells.Item(4, 49).Value = 3.75  # AW4: 4 -> 3.75
$ws.Cells.Item(6, 7).Value = 1.48  # G6: 1.45 -> 1.48
$ws.Cells.Item(6, 9).Value = 7  # I6: 7.5 -> 7
$ws.Cells.Item(6, 11).Value = 2.2  # K6: 2.25 -> 2.2
$ws.Cells.Item(6, 13).Value = 1.07  # M6: 1.06 -> 1.07
$ws.Cells.Item(6, 14).Value = 9  # N6: 9.5 -> 9
$ws.Cells.Item(6, 17).Value = 2.05  # Q6: 2.03 -> 2.05
$ws.Cells.Item(6, 18).Value = 1.8  # R6: 1.83 -> 1.8
$ws.Cells.Item(6, 25).Value = 8.5  # Y6: 9 -> 8.5
$ws.Cells.Item(6, 29).Value = 9  # AC6: 9.5 -> 9
$ws.Cells.Item(6, 30).Value = 7.5  # AD6: 8 -> 7.5
$ws.Cells.Item(6, 32).Value = 67  # AF6: 81 -> 67
$ws.Cells.Item(6, 34).Value = 17  # AH6: 15 -> 17
$ws.Cells.Item(6, 35).Value = 41  # AI6: 34 -> 41
$ws.Cells.Item(6, 36).Value = 23  # AJ6: 21 -> 23
$ws.Cells.Item(6, 45).Value = 201  # AS6: 151 -> 201
$ws.Cells.Item(7, 25).Value = 10  # Y7: 9.5 -> 10
$ws.Cells.Item(7, 41).Value = 6  # AO7: 6.5 -> 6
$ws.Cells.Item(8, 8).Value = 4.33  # H8: 4.5 -> 4.33
$ws.Cells.Item(8, 11).Value = 2.38  # K8: 2.4 -> 2.38
$ws.Cells.Item(8, 15).Value = 1.25  # O8: 1.22 -> 1.25
$ws.Cells.Item(8, 16).Value = 3.75  # P8: 4 -> 3.75
$ws.Cells.Item(8, 17).Value = 1.83  # Q8: 1.8 -> 1.83
$ws.Cells.Item(8, 18).Value = 2.03  # R8: 2 -> 2.03
$ws.Cells.Item(8, 26).Value = 9  # Z8: 8.5 -> 9
$ws.Cells.Item(8, 34).Value = 19  # AH8: 21 -> 19
$ws.Cells.Item(8, 36).Value = 23  # AJ8: 26 -> 23
$ws.Cells.Item(8, 37).Value = 81  # AK8: 101 -> 81
$ws.Cells.Item(8, 49).Value = 8.5  # AW8: 9 -> 8.5
$ws.Cells.Item(9, 7).Value = 1.33  # G9: 1.3 -> 1.33
$ws.Cells.Item(9, 9).Value = 9.5  # I9: 11 -> 9.5
$ws.Cells.Item(9, 12).Value = 8.5  # L9: 9 -> 8.5
$ws.Cells.Item(9, 13).Value = 1.04  # M9: 1.05 -> 1.04
$ws.Cells.Item(9, 14).Value = 13  # N9: 11 -> 13
$ws.Cells.Item(9, 15).Value = 1.25  # O9: 1.29 -> 1.25
$ws.Cells.Item(9, 16).Value = 3.75  # P9: 3.5 -> 3.75
$ws.Cells.Item(9, 17).Value = 1.85  # Q9: 1.9 -> 1.85
$ws.Cells.Item(9, 18).Value = 2  # R9: 1.95 -> 2
$ws.Cells.Item(9, 19).Value = 1.36  # S9: 1.4 -> 1.36
$ws.Cells.Item(9, 20).Value = 3  # T9: 2.75 -> 3
$ws.Cells.Item(9, 21).Value = 2.2  # U9: 2.38 -> 2.2
$ws.Cells.Item(9, 22).Value = 1.62  # V9: 1.53 -> 1.62
$ws.Cells.Item(9, 23).Value = 6.5  # W9: 5.5 -> 6.5
$ws.Cells.Item(9, 24).Value = 6  # X9: 5.5 -> 6
$ws.Cells.Item(9, 26).Value = 8  # Z9: 7.5 -> 8
$ws.Cells.Item(9, 27).Value = 12  # AA9: 13 -> 12
$ws.Cells.Item(9, 29).Value = 10  # AC9: 9.5 -> 10
$ws.Cells.Item(9, 30).Value = 9  # AD9: 9.5 -> 9
$ws.Cells.Item(9, 31).Value = 23  # AE9: 26 -> 23
$ws.Cells.Item(9, 32).Value = 81  # AF9: 101 -> 81
$ws.Cells.Item(9, 34).Value = 19  # AH9: 21 -> 19
$ws.Cells.Item(9, 35).Value = 41  # AI9: 51 -> 41
$ws.Cells.Item(9, 36).Value = 26  # AJ9: 34 -> 26
$ws.Cells.Item(9, 37).Value = 101  # AK9: 151 -> 101
$ws.Cells.Item(9, 38).Value = 67  # AL9: 81 -> 67
$ws.Cells.Item(9, 39).Value = 67  # AM9: 81 -> 67
$ws.Cells.Item(9, 40).Value = 3.2  # AN9: 3.1 -> 3.2
$ws.Cells.Item(9, 44).Value = 41  # AR9: 51 -> 41
$ws.Cells.Item(9, 45).Value = 151  # AS9: 201 -> 151
$ws.Cells.Item(9, 46).Value = 3  # AT9: 2.75 -> 3
$ws.Cells.Item(9, 47).Value = 10  # AU9: 11 -> 10
$ws.Cells.Item(9, 48).Value = 67  # AV9: 81 -> 67
$ws.Cells.Item(9, 49).Value = 9.5  # AW9: 10 -> 9.5
$ws.Cells.Item(9, 50).Value = 41  # AX9: 51 -> 41
$ws.Cells.Item(9, 51).Value = 41  # AY9: 51 -> 41
$ws.Cells.Item(9, 52).Value = 201  # AZ9: 251 -> 201
$ws.Cells.Item(9, 53).Value = 251  # BA9: 301 -> 251
$ws.Cells.Item(10, 7).Value = 1.65  # G10: 1.73 -> 1.65
$ws.Cells.Item(10, 8).Value = 3.4  # H10: 3.3 -> 3.4
$ws.Cells.Item(10, 9).Value = 5.75  # I10: 5.25 -> 5.75
$ws.Cells.Item(10, 33).Value = 900  # AG10: 201 -> 900
$ws.Cells.Item(11, 17).Value = 2.05  # Q11: 2.08 -> 2.05
$ws.Cells.Item(11, 18).Value = 1.75  # R11: 1.73 -> 1.75
$ws.Cells.Item(11, 33).Value = 1000  # AG11: 201 -> 1000
$ws.Cells.Item(12, 7).Value = 2.67  # G12: 2.65 -> 2.67
$ws.Cells.Item(12, 8).Value = 2.55  # H12: 2.57 -> 2.55
$ws.Cells.Item(12, 10).Value = 3.35  # J12: 3.3 -> 3.35
$ws.Cells.Item(12, 11).Value = 1.8  # K12: 1.82 -> 1.8
$ws.Cells.Item(12, 12).Value = 3.75  # L12: 3.7 -> 3.75
$ws.Cells.Item(12, 14).Value = 4.9  # N12: 4.95 -> 4.9
$ws.Cells.Item(12, 16).Value = 2.25  # P12: 2.27 -> 2.25
$ws.Cells.Item(12, 17).Value = 2.45  # Q12: 2.42 -> 2.45
$ws.Cells.Item(12, 18).Value = 1.42  # R12: 1.44 -> 1.42
$ws.Cells.Item(12, 19).Value = 1.55  # S12: 1.53 -> 1.55
$ws.Cells.Item(12, 20).Value = 2.15  # T12: 2.2 -> 2.15
$ws.Cells.Item(12, 21).Value = 1.93  # U12: 1.91 -> 1.93
$ws.Cells.Item(12, 24).Value = 12.5  # X12: 12 -> 12.5
$ws.Cells.Item(12, 25).Value = 10.25  # Y12: 10 -> 10.25
$ws.Cells.Item(12, 26).Value = 35  # Z12: 32 -> 35
$ws.Cells.Item(12, 27).Value = 28  # AA12: 27 -> 28
$ws.Cells.Item(12, 28).Value = 45  # AB12: 40 -> 45
$ws.Cells.Item(12, 29).Value = 5.7  # AC12: 5.8 -> 5.7
$ws.Cells.Item(12, 30).Value = 5.1  # AD12: 5.2 -> 5.1
$ws.Cells.Item(12, 31).Value = 15.5  # AE12: 15 -> 15.5
$ws.Cells.Item(12, 34).Value = 7.1  # AH12: 7.3 -> 7.1
$ws.Cells.Item(12, 36).Value = 11.25  # AJ12: 11 -> 11.25
$ws.Cells.Item(12, 38).Value = 35  # AL12: 32 -> 35
$ws.Cells.Item(12, 41).Value = 15.5  # AO12: 15 -> 15.5
$ws.Cells.Item(12, 42).Value = 25  # AP12: 24 -> 25
$ws.Cells.Item(12, 46).Value = 2.15  # AT12: 2.18 -> 2.15
$ws.Cells.Item(12, 47).Value = 6.9  # AU12: 6.8 -> 6.9
$ws.Cells.Item(12, 51).Value = 27  # AY12: 26 -> 27
$ws.Cells.Item(12, 52).Value = 100  # AZ12: 90 -> 100
$ws.Cells.Item(12, 54).Value = 400  # BB12: 350 -> 400
